# fix min_pack in master barang
# Adds a new "Min Qty (Jumlah renteng dalam satu dus)" column (K) to the
# Sheet1 header row of the barang-import template, matching the style of
# the existing header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell, next to the existing "Diskon" header (J4)
$ws.Range("K4").Value = "Min Qty (Jumlah renteng dalam satu dus)"

# Reuse the same header formatting (fill/font) already used by J4 instead
# of defining a new style
$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Widen the new column enough to show the full header text
$ws.Columns.Item(11).ColumnWidth = 38.140625

# Match the on-screen view state left behind by the edit
$excel.ActiveWindow.Zoom = 85
$ws.Range("J2").Select()
